$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = "Thomas Debiasi"
$ws.Range("B33").Value = "Thomas Debiasi | MAI UNA GIOIA"
$ws.Range("C33").Value = "Daniele Dalbosco | SdrumALA"
$ws.Range("D33").Value = "FEDERICO NICOLODI | U.S. Guarna"
$ws.Range("E33").Value = "Andrea Conzatti | FC Savignano"
$ws.Range("F33").Value = "Giacomo Gasparini | MAI UNA GIOIA"
